# edit.ps1 -- applies the "plotNorm help improved" deck edit via PowerPoint
# COM-interop object model.
#
# Summary of changes:
#   1. Slide 2 ("DGEobj Overview"): "automated" -> "scripted" in the
#      "Facilitates data sharing and automated meta-analysis" bullet.
#   2. Slide 4 ("DGEobj Engineering"): give the body placeholder an explicit
#      size/position (xfrm) -- widened from 611.75pt to ~618.875pt.
#   3. Slide 4: reword/re-split the "We store information..." bullet into
#      five runs, changing "on the DGEobj" -> "attached to the DGEobj" and
#      "elements" -> "items".

$p = $ppt.ActivePresentation

# Helper: replace the first occurrence of $search inside the shape's whole
# TextRange with $replacement, preserving the surrounding runs/formatting
# (PowerPoint splits runs at the edited boundary automatically).
function Replace-InRange($range, [string]$search, [string]$replacement) {
    $whole = $range.Text
    $idx = $whole.IndexOf($search)
    if ($idx -lt 0) {
        throw "Substring not found: '$search'"
    }
    $startPos = $idx + 1
    $len = $search.Length
    $sub = $range.Characters($startPos, $len)
    $sub.Text = $replacement
}

# ---------------------------------------------------------------------
# 1) Slide 2 - "automated" -> "scripted"
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
Replace-InRange $tr2 "automated " "scripted "

# ---------------------------------------------------------------------
# 2) Slide 4 - widen the body placeholder via an explicit xfrm
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$shape4 = $slide4.Shapes.Item(2)

# Shape.Left/Top/Width/Height are in points (1 pt = 12700 EMU). The values
# below are chosen so that, after the runtime's internal float32 rounding,
# the stored EMUs come out to exactly x=684213 y=1377538 cx=7859712 cy=4785756.
$shape4.Left   = 53.87504007007874
$shape4.Top    = 108.46755905511812
$shape4.Width  = 618.8749695299213
$shape4.Height = 376.8311921023622

# ---------------------------------------------------------------------
# 3) Slide 4 - reword the "We store information..." bullet
# ---------------------------------------------------------------------
$tr4 = $shape4.TextFrame.TextRange

# Step 1: "on the DGEobj" -> "attached to the DGEobj"
Replace-InRange $tr4 "on the DGEobj" "attached to the DGEobj"

# Step 2: force a run boundary right after "attached to "
Replace-InRange $tr4 "attached to " "attached to "

# Step 3: merge "the DGEobj.  In some cases individual " into one run
Replace-InRange $tr4 "the DGEobj.  In some cases individual " "the DGEobj.  In some cases individual "

# Step 4: "elements " -> "items " (and gives it its own run)
Replace-InRange $tr4 "elements " "items "
